$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8
$ws.Range("G8").Value = 2.42
$ws.Range("H8").Value = 3.15
$ws.Range("J8").Value = 1.08
$ws.Range("K8").Value = 6.4
$ws.Range("L8").Value = 1.39
$ws.Range("M8").Value = 2.77
$ws.Range("N8").Value = 2.15
$ws.Range("O8").Value = 1.62
$ws.Range("P8").Value = 1.5
$ws.Range("Q8").Value = 2.42
$ws.Range("R8").Value = 1.87
$ws.Range("S8").Value = 1.83
$ws.Range("T8").Value = 7
$ws.Range("V8").Value = 9.5
$ws.Range("W8").Value = 25
$ws.Range("X8").Value = 22
$ws.Range("Y8").Value = 35
$ws.Range("Z8").Value = 6.4
$ws.Range("AA8").Value = 6
$ws.Range("AB8").Value = 15
$ws.Range("AC8").Value = 80
$ws.Range("AD8").Value = 700
$ws.Range("AE8").Value = 8
$ws.Range("AI8").Value = 26
$ws.Range("AJ8").Value = 37

# Row 9
$ws.Range("G9").Value = 1.48
$ws.Range("H9").Value = 4.33
$ws.Range("I9").Value = 6.5
$ws.Range("J9").Value = 1.06
$ws.Range("K9").Value = 10
$ws.Range("N9").Value = 1.93
$ws.Range("O9").Value = 1.88
$ws.Range("R9").Value = 2.1
$ws.Range("S9").Value = 1.67
$ws.Range("T9").Value = 6.5
$ws.Range("U9").Value = 6.5
$ws.Range("V9").Value = 8.5
$ws.Range("W9").Value = 9.5
$ws.Range("Y9").Value = 29
$ws.Range("AA9").Value = 8.5
$ws.Range("AB9").Value = 21
$ws.Range("AC9").Value = 67
$ws.Range("AE9").Value = 15
$ws.Range("AF9").Value = 34
$ws.Range("AG9").Value = 21
$ws.Range("AH9").Value = 81
$ws.Range("AJ9").Value = 51

# Row 13
$ws.Range("L13").Value = 1.25
$ws.Range("M13").Value = 3.75

# Row 14
$ws.Range("J14").Value = 1.06
$ws.Range("K14").Value = 10

# Row 17
$ws.Range("G17").Value = 2.3
$ws.Range("I17").Value = 2.8
$ws.Range("L17").Value = 1.14
$ws.Range("M17").Value = 5.5
$ws.Range("V17").Value = 9.5
$ws.Range("AC17").Value = 29
$ws.Range("AJ17").Value = 23

# Row 21
$ws.Range("T21").Value = 14
$ws.Range("U21").Value = 14.5
$ws.Range("W21").Value = 21
$ws.Range("Y21").Value = 15.5
$ws.Range("AE21").Value = 18.5
$ws.Range("AF21").Value = 25
$ws.Range("AJ21").Value = 22

# Row 22
$ws.Range("L22").Value = 1.15
$ws.Range("M22").Value = 4.8
$ws.Range("N22").Value = 1.47
$ws.Range("O22").Value = 2.52
$ws.Range("P22").Value = 1.27
$ws.Range("Q22").Value = 3.4
$ws.Range("S22").Value = 2.65
$ws.Range("T22").Value = 16.5
$ws.Range("U22").Value = 24
$ws.Range("AA22").Value = 7.9
$ws.Range("AB22").Value = 11
$ws.Range("AC22").Value = 32
$ws.Range("AE22").Value = 11.5
$ws.Range("AF22").Value = 12.5
$ws.Range("AJ22").Value = 17.5

# Row 23
$ws.Range("O23").Value = 1.54

# Row 24
$ws.Range("J24").Value = 1.07
$ws.Range("K24").Value = 9
$ws.Range("N24").Value = 2.15
$ws.Range("O24").Value = 1.63

# Row 25
$ws.Range("L25").Value = 1.13
$ws.Range("M25").Value = 6

# Row 31
$ws.Range("G31").Value = 2.2
$ws.Range("H31").Value = 3.2
$ws.Range("I31").Value = 3.4
$ws.Range("U31").Value = 11
$ws.Range("W31").Value = 21
$ws.Range("X31").Value = 17
$ws.Range("AG31").Value = 12
$ws.Range("AH31").Value = 34
$ws.Range("AI31").Value = 26
$ws.Range("AJ31").Value = 34

# Row 35
$ws.Range("L35").Value = 1.22
$ws.Range("M35").Value = 4
$ws.Range("N35").Value = 1.75
$ws.Range("O35").Value = 2.05

# Row 36
$ws.Range("G36").Value = 2.7
$ws.Range("I36").Value = 2.25
$ws.Range("L36").Value = 1.33
$ws.Range("M36").Value = 3.25
$ws.Range("N36").Value = 2.03
$ws.Range("O36").Value = 1.78
$ws.Range("R36").Value = 1.83
$ws.Range("S36").Value = 1.83
$ws.Range("V36").Value = 11
$ws.Range("AE36").Value = 7.5
$ws.Range("AH36").Value = 21

# Row 42
$ws.Range("G42").Value = 2.6
$ws.Range("I42").Value = 2.6
$ws.Range("R42").Value = 1.7
$ws.Range("S42").Value = 2.05

# Row 45
$ws.Range("G45").Value = 2.1
$ws.Range("J45").Value = 1.04
$ws.Range("K45").Value = 13
$ws.Range("N45").Value = 1.9
$ws.Range("O45").Value = 1.9
$ws.Range("R45").Value = 1.75
$ws.Range("S45").Value = 2
$ws.Range("U45").Value = 10
$ws.Range("W45").Value = 19
$ws.Range("Z45").Value = 12

# Row 48
$ws.Range("G48").Value = 1.34
$ws.Range("H48").Value = 5.1
$ws.Range("I48").Value = 7.2
$ws.Range("O48").Value = 2.75
$ws.Range("P48").Value = 1.21
$ws.Range("Q48").Value = 3.95
$ws.Range("R48").Value = 1.65
$ws.Range("S48").Value = 2.12
$ws.Range("T48").Value = 10.5
$ws.Range("AA48").Value = 10.75
$ws.Range("AE48").Value = 28
$ws.Range("AH48").Value = 175
$ws.Range("AI48").Value = 65

# Row 52
$ws.Range("G52").Value = 1.91
$ws.Range("H52").Value = 3.3
$ws.Range("I52").Value = 3.8
$ws.Range("J52").Value = 1.08
$ws.Range("K52").Value = 8
$ws.Range("U52").Value = 8
$ws.Range("AF52").Value = 19
$ws.Range("AG52").Value = 15
$ws.Range("AJ52").Value = 51
